$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression
$ws.Range("B2").Value = 2971975655329742
$ws.Range("C2").Value = 2971975655329742
$ws.Range("D2").Value = 2971975655329742

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 16965079690292.61
$ws.Range("C3").Value = 5884185865701.873
$ws.Range("D3").Value = 105979022492249.6

# Row 4 - rename model and update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 4812777585615.108
$ws.Range("C4").Value = 5811478685486.691
$ws.Range("D4").Value = 21583745515484.32

# Row 5 - rename model and update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 151257572606756.2
$ws.Range("C5").Value = 614780850147377.9
$ws.Range("D5").Value = 1957504402503985
